$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

# Row 21: Diferencia Stock (L) goes from 2 to 0; uds. Objetivo semana pasada (R) goes from 0 to 2
$ws.Range("L21").Value = 0
$ws.Range("R21").Value = 2

# Row 26: uds. Objetivo semana pasada (R) goes from 0 to 1
$ws.Range("R26").Value = 1

# Row 28: uds. Objetivo semana pasada (R) goes from 0 to 2
$ws.Range("R28").Value = 2

# Row 30: uds. Objetivo semana pasada (R) goes from 0 to 2
$ws.Range("R30").Value = 2

# Row 31: uds. Objetivo semana pasada (R) 0->3, Tendencia Consumo (T) 9->6, Pedido Final (U) 11->8
$ws.Range("R31").Value = 3
$ws.Range("T31").Value = 6
$ws.Range("U31").Value = 8

# Totals: Total_Unidades (C36) 34->31, Total_Ajuste_Stock (C47) 2->0
$ws.Range("C36").Value = 31
$ws.Range("C47").Value = 0
